$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header C1 text from "Leer mas" to "Link"
$ws.Range("C1").Value = 'Link'

# Add new header D1 "Parrafos" with same style as other header cells (copy format from C1)
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("D1").Value = 'Párrafos'

# Populate column D (D2:D13) with paragraph lists extracted via NLP
$ws.Range("D2").Value = '[''21 de abril, 2025.'', ''¡La Universidad Manuela Beltrán sigue celebrando sus 50 años de historia todos los meses del año 2025 con invitados especiales.'', ''En este mes de abril , tenemos el honor de contar como país invitado a Finlandia 🇫🇮, en una jornada cultural y académica que nos lleva a explorar nuevas perspectivas y experiencias del norte de Europa.'', ''Conferencias, cine, foros, exposiciones visuales y vivencias internacionales llenan de conocimiento e inspiración a toda la comunidad Manuelista.'', ''¡Universidad Manuela Beltrán, 50 años transformando vidas!'']'
$ws.Range("D3").Value = '[''17 de febrero, 2024.'', ''Muestra toda tu creatividad y originalidad en este concurso que premia la innovación de los estudiantes que cursan entre primero y quinto semestre de los programas de pregrado de Bogotá, Cajicá, Bucaramanga y UMB Virtual.'', ''Debes presentar un invento que funcione, es decir, que opere, pero que no tenga ninguna utilidad, en otras palabras, que no sirva para nada. El ganador del primer puesto se llevará una media beca.'', ''Cada ítem se validará por medio del cumplimiento o no de las afirmaciones que describen cada aspecto, lo anterior asignará un puntaje de 10 si se cumple la afirmación y de 0 si la afirmación es falsa, estos puntajes se promediarán por cada aspecto, para posteriormente obtener la calificación final otorgando a cada ítem una ponderación del 20%.'', ''Los finalistas serán quienes tengan las 10 puntuaciones más altas por sede. Las afirmaciones que evaluarán los inventos se encuentran a continuación:'']'
$ws.Range("D4").Value = '[''14 de febrero, 2025.'', ''La Universidad Manuela Beltrán inició la celebración de sus 50 años de historia, un camino marcado por la excelencia, la innovación y el compromiso con la educación. Nuestra primera Jornada Académica y Cultural fue un encuentro inolvidable, en el que tuvimos el honor de contar con la participación de España como país invitado de honor.'', ''Desde el primer día, la emoción y el conocimiento se unieron en actividades académicas que nos transportaron más allá de nuestras fronteras: charlas inspiradoras sobre oportunidades de estudio en España, un fascinante recorrido por la magia del teatro español.'', ''La música también hizo vibrar nuestros corazones con la presentación del Coro del Centro Cultural y Educativo Español Reyes Católicos, que nos deleitó con composiciones contemporáneas, reflejando la riqueza y diversidad cultural de España. Así mismo, la exposición: Arte y Cultura a Través de las Letras, piezas de los ilustradores y diseñadores gráficos españoles, donde se establece una conexión entre imagen y lenguaje.'', ''Así mismo se desarrolló el Cine Foro, actividad donde España presentó películas enmarcadas en Derechos Humanos. Frágil Equilibrio, La Boca del Lobo, Acorazado, entre otras, formaron parte de este espacio de cine.'', ''La antigüedad Grecorromana y la cultura LGBTI contemporánea: memorias, imaginarios y referentes de resistencia, Acción colectiva y cambio climático: El activismo ambiental en perspectiva y Psicología política y derechos humanos, fueron algunas de las charlas que también se brindaron en este espacio académico.'', ''Este 2025 celebramos medio siglo de sueños cumplidos, de transformación y de impacto en miles de estudiantes.'', ''¡Universidad Manuela Beltrán, 50 años transformando vidas!'']'
$ws.Range("D5").Value = '[''13 de febrero, 2025.'', ''Frágil Equilibrio, La Boca del Lobo, Acorazado, entre otras, formaron parte de esta actividad desarrollada para la Comunidad Manuelista.'']'
$ws.Range("D6").Value = '[''13 de febrero, 2025.'', ''Otra de las actividades desarrolladas en el marco de la celebración de los 50 años de la Universidad con España como país invitado, fue la que se realizó el pasado 12 de febrero denominada: La antigüedad Grecorromana y la cultura LGBTI contemporánea: memorias, imaginarios y referentes de resistencia.'', ''Este evento académico, realizado en el auditorio de Cine y Televisión, fue desarrollarlo por Ricardo del Molino, Doctor y Licenciado en Humanidades de la Universidad Carlos lll y Licenciado en Derecho de la Universidad Autónoma Madrid.'', ''Este espacio académico contó con la gran participación de los estudiantes de Psicología y Cine y Televisión.'']'
$ws.Range("D7").Value = '[''16 de enero, 2025.'', ''Durante la primera Semana de Inducción 2025-1, nuestros estudiosos tuvieron la oportunidad de conocer nuestras modernas instalaciones, interactuar con docentes, y sumergirse en la vida universitaria Manuelista'', ''Charlas, talleres, capacitaciones, conciertos y mucho más fueron las actividades que se desarrollaron durante esta primera semana de inducción.'', ''¡Les deseamos un exitoso inicio de semestre!'']'
$ws.Range("D8").Value = '[''26 de noviembre, 2024.'', ''En esta oportunidad, el ganador fue el estudioso Davis Esteven Nausan Martínez del programa de Ingeniería Biomédica de cuarto semestre con su proyecto Gelease. Su invento consistió en crear un dispositivo para evitar el contacto del gel con las manos, cuando una persona se realice un peinado. Davis Nausan ganó una media beca para estudiar el próximo semestre.'', ''¡Felicitaciones a todos los estudiantes que participaron este 2025 – 2 con sus grandes inventos! 😄'']'
$ws.Range("D9").Value = '[''6 de noviembre, 2024.'']'
$ws.Range("D10").Value = '[''6 de noviembre, 2024.'', ''Nuestra selección de Fútbol Tenis se coronó campeona en diversas categorías del Grupo Cerros, Torneo Deportivo Universitario. Nuestros Manuelistas demostraron un juego excepcional.'', ''Así mismo logramos el primer lugar, en la categoría #3 de Tenis en este importante torneo deportivo, que reúne a las principales Instituciones de Educación Superior.'', ''¡Estos triunfos es el resultado de mucho esfuerzo y dedicación!'', ''Santiago Santander – Fisioterapia'', ''Santiago Corredor – Fisioterapia'']'
$ws.Range("D11").Value = '[''28 de octubre, 2024.'', ''Se desarrolló en la Universidad el Foro Internacional FFC es espacio creado por Fútbol con Corazón con el objetivo de visibilizar las apuestas en diferentes temáticas que impactan comunidades a nivel mundial. En su octava versión, el foro se enfoca en socializar alternativas innovadoras de medición de impacto y buenas prácticas globales en monitoreo y evaluación de proyectos sociales.'', ''En este evento, en el que Bogotá se convierte en el epicentro del deporte para el desarrollo social en América Latina, participaron diferentes personalidades y figuras representativas del deporte como:'', ''La Universidad Manuela Beltrán se une a esta gran espacio de espacio de diálogo transformador, donde se reunieron los conocimientos y experiencias de académicos, deportistas, líderes comunitarios y representantes de los sectores público y privado.'']'
$ws.Range("D12").Value = '[''15 de octubre, 2024.'']'
$ws.Range("D13").Value = '[''2 de octubre, 2024'', ''La Universidad Manuela Beltrán celebra con orgullo los 40 AÑOS de este evento, un hito que marca cuatro décadas de compromiso con la inclusión y el deporte. Este 1 de octubre se desarrolló la inauguración en el Gimnasio Vermont donde más de 1.000 deportistas de 50 fundaciones de Bogotá y Cundinamarca se reunieron para demostrar su talento y espíritu competitivo.'', ''Durante estos 40 años, las Olimpiadas han impactado la vida de más de 57,000 personas con discapacidad física, mental y sensorial. A través de la recreación y el deporte, hemos mejorado la calidad de vida de estos participantes y sus familias, demostrando que la inclusión social es posible. Niños, niñas, adultos y personas mayores han formado parte durante todo este tiempo, gracias a las fundaciones que año tras año participan en este evento deportivo y cultural.'', ''¡Todos Somos Campeones!'']'
